# Fix state view colours and school attendance status.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("Data")
$wsDesc = $wb.Worksheets.Item("Description")

# Update the school attendance status text (Description sheet, row "Status").
$wsDesc.Range("B3").Value = "Not on track"

# Update selection/active-cell state on the "Data" sheet and make it inactive.
$null = $wsData.Select()
$null = $wsData.Range("A32").Select()

# Update selection/active-cell state on the "Description" sheet and make it
# the active tab (matches activeTab="1" in workbook.xml bookViews).
$null = $wsDesc.Select()
$null = $wsDesc.Range("B5").Select()
